# Apply the diff: update lab-result labels, units, and reference ranges to
# normalize dash/tilde/caret typography, fix a few mislabeled rows, and
# clear some stray single-character cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericText {
    param(
        [string]$Addr,
        [string]$Val
    )
    $rng = $ws.Range($Addr)
    # These values look like plain numbers (e.g. "4.9") - force text storage
    # so they aren't silently reinterpreted as numeric cells.
    $rng.NumberFormat = "@"
    $rng.Value = $Val
}

# Row 2 - WBC
$ws.Range("B2").Value = "中性细胞数"
Set-NumericText "C2" "4.9"
$ws.Range("E2").Value = "4-10"

# Row 3 - NEUT#
$ws.Range("B3").Value = "中性粒细胞数"
$ws.Range("E3").Value = "2-7.5"

# Row 4 - LYMPH#
$ws.Range("B4").Value = "淋巴细胞数"
Set-NumericText "C4" "1.3"
$ws.Range("D4").Value = "10^9/L"
$ws.Range("E4").Value = "0.8-4"

# Row 5 - MONO#
$ws.Range("B5").Value = "单核细胞数"
$ws.Range("D5").Value = "10^9/L"
$ws.Range("E5").Value = "0.16-1.2"

# Row 6 - E0#
$ws.Range("B6").Value = "嗜酸性粒细胞数"
$ws.Range("D6").Value = "10^9/L"
$ws.Range("E6").Value = "0.02-0.5"

# Row 7 - BASO#
$ws.Range("B7").Value = "嗜碱性粒细胞"
$ws.Range("D7").Value = "10^9/L"
$ws.Range("E7").Value = "0-0.1"

# Row 8 - NEUT%
$ws.Range("B8").Value = "中性粒细胞数"
$ws.Range("E8").Value = "50-75"

# Row 9 - LYMPH%
$ws.Range("B9").Value = "淋巴细胞比率"
Set-NumericText "C9" "26.5"
$ws.Range("E9").Value = "20-40"

# Row 10 - MONO%
$ws.Range("B10").Value = "单核细胞数"
Set-NumericText "C10" "5.1"
$ws.Range("E10").Value = "4-12"

# Row 11 - E0%
$ws.Range("B11").Value = "嗜酸性粒细胞数"
$ws.Range("E11").Value = "0.5-5"

# Row 12 - BAS0%
$ws.Range("B12").Value = "嗜碱性粒细胞"
$ws.Range("E12").Value = "0-1"

# Row 13 - RBC
$ws.Range("B13").Value = "中性细胞数"
$ws.Range("D13").Value = "10^12/L"
$ws.Range("E13").Value = "3.5-5.5"

# Row 14 - HGB
$ws.Range("E14").Value = "110-160"

# Row 15 - HCT
$ws.Range("E15").Value = "37-49"

# Row 16 - MCV
$ws.Range("D16").Value = "fL"
$ws.Range("E16").Value = "82-95"

# Row 17 - MCH
$ws.Range("B17").Value = "平均血红蛋白量"
$ws.Range("E17").Value = "27-31"

# Row 18 - MCHC
$ws.Range("B18").Value = "平均血红蛋白浓度"
$ws.Range("E18").Value = "320-360"

# Row 19 - RDW-SD
$ws.Range("B19").Value = "红细胞分布宽度"
$ws.Range("D19").Value = "fL"
$ws.Range("E19").Value = "37-54"

# Row 20 - RDW-CV
$ws.Range("B20").Value = "红细胞分布宽度"
$ws.Range("E20").Value = "11-16"

# Row 21 - PLT
$ws.Range("E21").Value = "100-300"

# Row 22 - PDW
$ws.Range("D22").Value = "fL"
$ws.Range("E22").Value = "9-17"

# Row 23 - MPV
$ws.Range("D23").Value = "fL"
$ws.Range("E23").Value = "9.4-12.5"

# Row 24 - PCT
$ws.Range("E24").Value = "0.11-0.27"

# Row 25 - P-LCR
$ws.Range("E25").Value = "13-43"

# Row 26
$ws.Range("B26").Value = "AST/ALT"

# Row 27 - clear the stray "1"
$ws.Range("B27").ClearContents()

# Row 31 - clear the stray "！！"
$ws.Range("B31").ClearContents()

# Row 32 - clear the stray "！"
$ws.Range("B32").ClearContents()
